$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Correct the auto date placeholder text ("11/2/2023" -> "12/21/2023")
#    on the slide master and every slide layout (ppPlaceholderDate = 16).
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($container, $newDateText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }

        $isPlaceholder = $false
        $placeholderType = -1
        try {
            $placeholderType = $shp.PlaceholderFormat.Type
            $isPlaceholder = $true
        } catch {
            $isPlaceholder = $false
        }

        if ($isPlaceholder -and $placeholderType -eq 16) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -ne $newDateText) {
                $tr.Text = $newDateText
            }
        }
    }
}

$newDateText = "12/21/2023"

$master = $p.SlideMaster
Update-DatePlaceholders $master $newDateText

$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholders $layouts.Item($j) $newDateText
}

# ---------------------------------------------------------------------------
# 2) Fix the timing callout on the slide: "60min" -> "8" + "0min"
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 12" -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "60min") {
            $tr.Text = "8"
            [void]$tr.InsertAfter("0min")
        }
    }
}
